$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New band: The Rolling Stones, Band ID 781 - three performance rows (10-12)
$ws.Range("A10").Value = 781
$ws.Range("B10").Value = 8437
$ws.Range("A11").Value = 781
$ws.Range("B11").Value = 4372
$ws.Range("A12").Value = 781
$ws.Range("B12").Value = 1834
$ws.Rows.Item(10).RowHeight = 15

# Label cell (merged D10:E12), matching the style used for the other band labels
$ws.Range("D7:E8").Copy()
$ws.Range("D10:E12").PasteSpecial(-4122)
$ws.Range("D10").Value = "#The Rolling Stones´s Band ID: 781"
$ws.Range("D10:E12").Merge()

# SQL-insert formula, following the same pattern as the rest of the column
$ws.Range("G10").Formula = "=CONCATENATE(`$K`$1,A10,`$K`$3,B10,`$K`$2)"
$ws.Range("G11").Formula = "=CONCATENATE(`$K`$1,A11,`$K`$3,B11,`$K`$2)"
$ws.Range("G12").Formula = "=CONCATENATE(`$K`$1,A12,`$K`$3,B12,`$K`$2)"

# Reflect the new selection state
$ws.Range("G10:G12").Select()
